$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SLS_dict")

$ws.Range("B4").Value = "Bei **wahren** Aussagen drücken Sie bitte die **J-Taste** (für „ja, wahr"").\\`n Bei **falschen** Aussagen drücken Sie bitte die **F-Taste** (für „falsch"") .\\**Bitte legen Sie Ihre Zeigefinger jetzt auf die F- und J-Taste**.\\ "
$ws.Range("B16").Value = "Weiter mit F oder J"
$ws.Range("C16").Value = "Continue with F or J"
$ws.Range("C4").Value = "For **true** statements, please press the **J-key** .\\`n For **false** statements, please press the **F-key** .\\You can best use your index fingers for the F and J keys.\\"

$ws.Range("C4").Select()
